# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates column G ("K") values on Sheet1 for rows 2-24 (row 22 unchanged at 0)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 0
    3  = 3
    4  = 1
    5  = 1
    6  = 0
    7  = 2
    8  = 0
    9  = 2
    10 = 1
    11 = 2
    12 = 2
    13 = 2
    14 = 3
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 3
    20 = 1
    21 = 2
    23 = 0
    24 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
